# chore: update content seeder
#
# - "content" sheet: drop the image/title/description columns in favour of
#   ustadzName/bunnyId/url, and swap the sample rows to the new
#   Ustadz-Hanan-Attaki / bunny-id seed data.
# - "contentCategory" sheet keeps its VLOOKUP formulas pointed at
#   content!A:D, which now resolves to #REF! because content only has
#   3 data columns after the drop (left as-is, matching upstream).
# - Active tab moves from contentCategory to content.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "content" sheet: remove the old "image" column. Excel shifts every
#    column after it left by one and auto-rewrites any formula that
#    referenced the old full-column range (contentCategory!B2:B4 goes
#    from content!A:D to content!A:C here), which is exactly what the
#    target workbook shows (VLOOKUP now errors with #REF!).
# ---------------------------------------------------------------------
$content = $wb.Worksheets.Item("content")
$content.Columns.Item(3).Delete()

# Columns are now: A id | B type | C title | D description | E link |
#                  F price | G createBy | H isActive
# Re-purpose C/D/E as the new ustadzName/bunnyId/url columns.

$content.Range("C1").Value = "ustadzName"
$content.Range("D1").Value = "bunnyId"
$content.Range("E1").Value = "url"

$content.Range("C2").Value = "Ustadz Hanan Attaki"
$content.Range("D2").Value = "132a4sd6f8as7d9g"
$content.Range("E2").Value = "https://www.youtube.com/watch?v=123456"

$content.Range("C3").Value = "Ustadz Hanan Attaki"
$content.Range("D3").Value = "132a4sd6f8as7d9f"
$content.Range("E3").Value = "https://www.youtube.com/watch?v=123456"

# ---------------------------------------------------------------------
# 2. View state: "content" becomes the selected/active sheet (was
#    "contentCategory"); give it the E2 selection the source file shows.
# ---------------------------------------------------------------------
$content.Activate()
$content.Range("E2").Select()

# ---------------------------------------------------------------------
# 3. "contentCategory" sheet keeps its data, but its selection/scroll
#    position changes and it's no longer the tab shown on open.
# ---------------------------------------------------------------------
$contentCategory = $wb.Worksheets.Item("contentCategory")
$contentCategory.Activate()
$contentCategory.Range("C21:C22").Select()

# Put the focus back on "content" last so it ends up the active tab.
$content.Activate()
